{"js": "// Update the license name from \"CC BY-NC-SA 4.0 ES\" to\n// \"CC BY-NC-SA 4.0 International\" while preserving the existing\n// (bold) character formatting of the run that contains it.\n\nconst results = context.document.body.search(\"CC BY-NC-SA 4.0 ES\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"CC BY-NC-SA 4.0 International\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the license name from \"CC BY-NC-SA 4.0 ES\" to\n# \"CC BY-NC-SA 4.0 International\" while preserving the existing\n# (bold) character formatting of the run that contains it.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"CC BY-NC-SA 4.0 ES\"\n$find.Replacement.Text = \"CC BY-NC-SA 4.0 International\"\n$find.Forward = $true\n$find.Wrap = $wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $find.MatchSoundsLike,\n    $find.MatchAllWordForms,\n    $find.Forward,\n    $find.Wrap,\n    $find.Format,\n    $find.Replacement.Text,\n    $wdReplaceOne\n)\n"}
